# Auto-generated edit script updating profit/price figures across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 68000
$ws.Range("J3").Value = 68000
$ws.Range("L3").Value = 68000
$ws.Range("N3").Value = -68228
$ws.Range("H102").Value = 68000
$ws.Range("J102").Value = 68000
$ws.Range("L102").Value = 68000
$ws.Range("N102").Value = -74490
$ws.Range("H103").Value = 504.7143
$ws.Range("I103").Value = 483.5
$ws.Range("J103").Value = 513.2
$ws.Range("K103").Value = 1450.5
$ws.Range("L103").Value = 1539.6
$ws.Range("M103").Value = -864.5
$ws.Range("N103").Value = -2711.6
$ws.Range("H106").Value = 718838.2
$ws.Range("I106").Value = 837394.7
$ws.Range("J106").Value = 7499.5
$ws.Range("K106").Value = 837394.7
$ws.Range("L106").Value = 7499.5
$ws.Range("M106").Value = -836763.7
$ws.Range("N106").Value = -8761.5
$ws.Range("H112").Value = 3425.3572
$ws.Range("J112").Value = 3381.1538
$ws.Range("L112").Value = 10143.4614
$ws.Range("N112").Value = -12359.4614
$ws.Range("H129").Value = 1510.4
$ws.Range("I129").Value = 720.75
$ws.Range("J129").Value = 2036.8334
$ws.Range("K129").Value = 2162.25
$ws.Range("L129").Value = 6110.5002
$ws.Range("M129").Value = 2837.75
$ws.Range("N129").Value = -16110.5002
$ws.Range("H137").Value = 3817.8096
$ws.Range("I137").Value = 2655.7812
$ws.Range("J137").Value = 7536.3
$ws.Range("K137").Value = 7967.3436
$ws.Range("L137").Value = 22608.9
$ws.Range("M137").Value = -5417.3436
$ws.Range("N137").Value = -27708.9
$ws.Range("H138").Value = 5225.2
$ws.Range("J138").Value = 6215.8696
$ws.Range("L138").Value = 18647.6088
$ws.Range("N138").Value = -28927.6088

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 18071
$ws.Range("I2").Value = 18833
$ws.Range("J2").Value = 17499.5
$ws.Range("K2").Value = 18833
$ws.Range("L2").Value = 17499.5
$ws.Range("M2").Value = -18720
$ws.Range("N2").Value = -17725.5
$ws.Range("H32").Value = 701.27
$ws.Range("I32").Value = 585.05206
$ws.Range("J32").Value = 3490.5
$ws.Range("K32").Value = 585.05206
$ws.Range("L32").Value = 3490.5
$ws.Range("M32").Value = -298.05206
$ws.Range("N32").Value = -4064.5
$ws.Range("H45").Value = 1316.5238
$ws.Range("I45").Value = 1105.25
$ws.Range("J45").Value = 1992.6
$ws.Range("K45").Value = 1105.25
$ws.Range("L45").Value = 1992.6
$ws.Range("M45").Value = -728.25
$ws.Range("N45").Value = -2746.6
$ws.Range("H46").Value = 7899.25
$ws.Range("J46").Value = 8519.857
$ws.Range("L46").Value = 8519.857
$ws.Range("N46").Value = -9157.857
$ws.Range("H104").Value = 36741.668
$ws.Range("J104").Value = 36741.668
$ws.Range("L104").Value = 36741.668
$ws.Range("N104").Value = -43729.668
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H116").Value = 18071
$ws.Range("I116").Value = 18833
$ws.Range("J116").Value = 17499.5
$ws.Range("K116").Value = 18833
$ws.Range("L116").Value = 17499.5
$ws.Range("M116").Value = -16539
$ws.Range("N116").Value = -22087.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 18071
$ws.Range("I3").Value = 18833
$ws.Range("J3").Value = 17499.5
$ws.Range("K3").Value = 18833
$ws.Range("L3").Value = 17499.5
$ws.Range("M3").Value = -18719
$ws.Range("N3").Value = -17727.5
$ws.Range("H15").Value = 9000
$ws.Range("J15").Value = 9000
$ws.Range("L15").Value = 9000
$ws.Range("N15").Value = -9454
$ws.Range("H95").Value = 11414
$ws.Range("J95").Value = 11414
$ws.Range("L95").Value = 11414
$ws.Range("N95").Value = -16906

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4522.2856
$ws.Range("I31").Value = 1369
$ws.Range("J31").Value = 6015.9473
$ws.Range("K31").Value = 1369
$ws.Range("L31").Value = 6015.9473
$ws.Range("M31").Value = -1074
$ws.Range("N31").Value = -6605.9473
$ws.Range("H34").Value = 4522.2856
$ws.Range("I34").Value = 1369
$ws.Range("J34").Value = 6015.9473
$ws.Range("K34").Value = 1369
$ws.Range("L34").Value = 6015.9473
$ws.Range("M34").Value = -1167
$ws.Range("N34").Value = -6419.9473
$ws.Range("H58").Value = 5292.081
$ws.Range("I58").Value = 4753.9
$ws.Range("J58").Value = 5925.2354
$ws.Range("K58").Value = 4753.9
$ws.Range("L58").Value = 5925.2354
$ws.Range("M58").Value = -4550.9
$ws.Range("N58").Value = -6331.2354
$ws.Range("H62").Value = 8213.5
$ws.Range("I62").Value = 2962
$ws.Range("J62").Value = 10464.143
$ws.Range("K62").Value = 2962
$ws.Range("L62").Value = 10464.143
$ws.Range("N62").Value = -11712.143
$ws.Range("H65").Value = 8213.5
$ws.Range("I65").Value = 2962
$ws.Range("J65").Value = 10464.143
$ws.Range("K65").Value = 14810
$ws.Range("L65").Value = 52320.715
$ws.Range("N65").Value = -58560.715
$ws.Range("H96").Value = 17800
$ws.Range("J96").Value = 17800
$ws.Range("L96").Value = 17800
$ws.Range("N96").Value = -23292
$ws.Range("H99").Value = 5364.364
$ws.Range("I99").Value = 5527.5
$ws.Range("J99").Value = 4929.3335
$ws.Range("K99").Value = 5527.5
$ws.Range("L99").Value = 4929.3335
$ws.Range("M99").Value = -4029.5
$ws.Range("N99").Value = -7925.3335
$ws.Range("H122").Value = 1923.7715
$ws.Range("I122").Value = 1636.64
$ws.Range("K122").Value = 4909.92
$ws.Range("M122").Value = -2459.92
$ws.Range("H126").Value = 5364.364
$ws.Range("I126").Value = 5527.5
$ws.Range("J126").Value = 4929.3335
$ws.Range("K126").Value = 16582.5
$ws.Range("L126").Value = 14788.0005
$ws.Range("M126").Value = -14112.5
$ws.Range("N126").Value = -19728.0005
$ws.Range("H132").Value = 21341.488
$ws.Range("I132").Value = 26313.812
$ws.Range("K132").Value = 78941.436
$ws.Range("M132").Value = -76411.436
$ws.Range("H136").Value = 5292.081
$ws.Range("I136").Value = 4753.9
$ws.Range("J136").Value = 5925.2354
$ws.Range("K136").Value = 14261.7
$ws.Range("L136").Value = 17775.7062
$ws.Range("M136").Value = -11711.7
$ws.Range("N136").Value = -22875.7062

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 420.35294
$ws.Range("I14").Value = 420.35294
$ws.Range("K14").Value = 1261.05882
$ws.Range("M14").Value = -1088.05882
$ws.Range("H107").Value = 4348.524
$ws.Range("J107").Value = 4756.0703
$ws.Range("L107").Value = 14268.2109
$ws.Range("N107").Value = -18108.2109
$ws.Range("H122").Value = 116642.87
$ws.Range("I122").Value = 566.5
$ws.Range("J122").Value = 134970.7
$ws.Range("K122").Value = 5098.5
$ws.Range("L122").Value = 1214736.3
$ws.Range("M122").Value = -2648.5
$ws.Range("N122").Value = -1219636.3
$ws.Range("H127").Value = 3949.5
$ws.Range("J127").Value = 3949.5
$ws.Range("L127").Value = 11848.5
$ws.Range("N127").Value = -21768.5
$ws.Range("H131").Value = 41275332
$ws.Range("I131").Value = 48485324
$ws.Range("J131").Value = 33344344
$ws.Range("K131").Value = 145455972
$ws.Range("L131").Value = 100033032
$ws.Range("M131").Value = -145450932
$ws.Range("N131").Value = -100043112

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 18718.273
$ws.Range("J15").Value = 18718.273
$ws.Range("L15").Value = 18718.273
$ws.Range("N15").Value = -19294.273
$ws.Range("H49").Value = 29000
$ws.Range("I49").Value = 29000
$ws.Range("K49").Value = 29000
$ws.Range("H81").Value = 18718.273
$ws.Range("J81").Value = 18718.273
$ws.Range("L81").Value = 18718.273
$ws.Range("N81").Value = -20714.273
$ws.Range("H84").Value = 18718.273
$ws.Range("J84").Value = 18718.273
$ws.Range("L84").Value = 56154.819
$ws.Range("N84").Value = -66138.819
$ws.Range("H122").Value = 4386.9375
$ws.Range("I122").Value = 4071.08
$ws.Range("K122").Value = 12213.24
$ws.Range("M122").Value = -9763.24

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3881.6924
$ws.Range("I40").Value = 2353.8333
$ws.Range("J40").Value = 7319.375
$ws.Range("K40").Value = 2353.8333
$ws.Range("L40").Value = 7319.375
$ws.Range("M40").Value = -2217.8333
$ws.Range("N40").Value = -7591.375
$ws.Range("H120").Value = 33332.668
$ws.Range("J120").Value = 33332.668
$ws.Range("L120").Value = 33332.668
$ws.Range("N120").Value = -43008.668
$ws.Range("H132").Value = 4246.25
$ws.Range("I132").Value = 4000
$ws.Range("J132").Value = 4328.3335
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 12985.0005
$ws.Range("M132").Value = -9470
$ws.Range("N132").Value = -18045.0005
$ws.Range("H136").Value = 12294.6
$ws.Range("I136").Value = 14589.8
$ws.Range("J136").Value = 9999.4
$ws.Range("K136").Value = 43769.39999999999
$ws.Range("L136").Value = 29998.2
$ws.Range("M136").Value = -41219.39999999999
$ws.Range("N136").Value = -35098.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1391.8823
$ws.Range("I113").Value = 1391.8823
$ws.Range("K113").Value = 4175.6469
$ws.Range("M113").Value = -2005.6469
$ws.Range("H136").Value = 5212.45
$ws.Range("I136").Value = 4292.92
$ws.Range("J136").Value = 6745
$ws.Range("K136").Value = 12878.76
$ws.Range("L136").Value = 20235
$ws.Range("M136").Value = -10328.76
$ws.Range("N136").Value = -25335
